$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "11/03/2013",
    "14/05/2013",
    "17/07/2013",
    "19/09/2013",
    "22/11/2013",
    "25/01/2014",
    "30/03/2014",
    "02/06/2014",
    "05/08/2014",
    "08/10/2014",
    "11/12/2014",
    "13/02/2015",
    "18/04/2015",
    "21/06/2015",
    "24/08/2015",
    "27/10/2015",
    "30/12/2015",
    "03/03/2016",
    "06/05/2016",
    "09/07/2016",
    "14/11/2016",
    "15/03/2017"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$i]
}

$ws.Range("D2:D23").Select()
